# Update the threshold table on Sheet1 and leave the selection on C3,
# matching the authored edit (values for alpha_distance_range /
# beta_distance_range rows were refined, and the cursor ended up on C3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.8
$ws.Range("C2").Value = 12.9
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 11.7

$ws.Range("C3").Select()
